$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$passwords = @{
    2  = "2zfWuPq"
    3  = "H2keyrt"
    4  = "lMN4ZYH"
    5  = "utm6nz1"
    6  = "FifDwaR"
    7  = "Wx5kjLY"
    8  = "IkB6Lzf"
    9  = "T0orWdn"
    10 = "8aofEAE"
    11 = "K7NVKG5"
    12 = "45UUO3R"
    13 = "2uGatTG"
    14 = "8ObgPzL"
    15 = "U1BWtcV"
    16 = "XA7mfYE"
    17 = "drvOaIq"
    18 = "BKHFfcj"
    19 = "MqSKwt7"
    20 = "VuJQsyS"
    21 = "gGKCX71"
    22 = "m4xtMvd"
    23 = "gxx6ZeU"
    24 = "JeK5w0d"
    25 = "6ASes5c"
    26 = "03ucRkh"
    27 = "W2kVi52"
    28 = "mV8JY1B"
    29 = "AmeI6Ee"
    30 = "DKC2trB"
    31 = "qwlZOpn"
}

foreach ($row in $passwords.Keys) {
    $ws.Cells.Item($row, 4).Value = $passwords[$row]
}
